# Update cryptocurrency price/volume data per the "Sat Jan 14 05:15:01 UTC 2023" symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D = Price, E = Volume(1h); values are stored as text, so force the "@" (Text)
# number format before writing so Excel does not coerce the numeric-looking
# strings (and trailing zeros / "%" suffixes) into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '307.91'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '6.33%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '32.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.94%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.339'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.41%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07426'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '11.24%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.760'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.27%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.698'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '8.47%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.570'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '15.54%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9231'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.15%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.01654'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2,453.80%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1673'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.36%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07665'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '14.35%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07952'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.79%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03072'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '4.66%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09837'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '9.46%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001532'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.72%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04557'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.75%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006504'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.66%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.470'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.45%'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.19%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3266'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.67%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1326'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.24%'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.55%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1630'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '4.00%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001216'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.92%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004511'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.21%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001170'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-6.48%'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001741'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04510'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '6.70%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007420'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '9.94%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1364'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '9.93%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002260'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '14.05%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01381'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '13.94%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006134'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '8.22%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01300'
